# Swap the species-identification data between row 4 and row 5
# (columns A, B, D, E, F, G, H, I). All other columns already match
# between the two rows, so this reduces to writing the new literal
# values described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (becomes the former row-5 species) ---
$ws.Range("A4").Value = 131116964
$ws.Range("B4").Value = 56762
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 100092
$ws.Range("F4").Value = "Större brunfladdermus"
$ws.Range("G4").Value = "Nyctalus noctula"
$ws.Range("H4").Value = "(Schreber, 1774)"
$ws.Range("I4").Value = "'42"

# --- Row 5 (becomes the former row-4 species) ---
$ws.Range("A5").Value = 131117036
$ws.Range("B5").Value = 56748
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 205998
$ws.Range("F5").Value = "Nordfladdermus"
$ws.Range("G5").Value = "Eptesicus nilssonii"
$ws.Range("H5").Value = "(A.Keyserling & Blasius, 1839)"
$ws.Range("I5").Value = "'443"
